# Apply "1st changes of mifos to finflux" edit:
# On the "Repayment Schedule" sheet, insert a new (blank) column before
# column N (14), shifting the old N/O/P columns one place to the right
# (N->O, O->P, P->Q). This matches the XML diff where the "Late" header
# and the "Over Due" data that used to live in N/O now live in O/P, with
# a brand-new blank column N, and the sheet's used range grows from
# A1:P15 to A1:Q15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new column at N (14); existing columns N.. shift right to O..
$ws.Columns.Item(14).Insert()

# Excel carries the left-neighbour column's width onto a freshly inserted
# column; give the new (blank) N column the same display width as M.
$ws.Columns.Item(14).ColumnWidth = 10.25

# Move the selection to match the post-edit state captured in the diff.
$ws.Range("R7").Select() | Out-Null
